{"js": "// 1. Update the letter date from \"September 19, 2025\" to \"September 21, 2025\".\nconst body = context.document.body;\nconst dateResults = body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2. Split the mailing-address paragraph \"2940 Sanor Pl, Santa Clara CA 95051\"\n//    (the sender/recipient block near the top of the letter, not the\n//    \"PROPERTY ADDRESS:\" line further down in the details table) into two\n//    paragraphs: \"2940 Sanor Pl\" and a new paragraph \"Santa Clara, CA 95051\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"2940 Sanor Pl, Santa Clara CA 95051\") {\n    addressParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (addressParagraph) {\n  // Insert the new second line right after the existing paragraph; it\n  // inherits the same paragraph/run formatting (Arial, 22 half-points).\n  addressParagraph.insertParagraph(\"Santa Clara, CA 95051\", Word.InsertLocation.after);\n  // Trim the original paragraph down to just the street address.\n  addressParagraph.getRange().insertText(\"2940 Sanor Pl\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Remove the now-superfluous empty \"No Spacing\" paragraph that sat\n//    directly below \"...Board of Directors\".\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"text,style\");\nawait context.sync();\n\nlet boardIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    boardIndex = i;\n    break;\n  }\n}\n\nif (boardIndex !== -1 && boardIndex + 1 < paragraphs2.items.length) {\n  const candidate = paragraphs2.items[boardIndex + 1];\n  if (candidate.text === \"\" && candidate.style === \"No Spacing\") {\n    candidate.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date from \"September 19, 2025\" to \"September 21, 2025\".\n$dateRange = $d.Content\n$find = $dateRange.Find\n$find.ClearFormatting()\n$find.Text = \"September 19, 2025\"\n$found = $find.Execute()\nif ($found) {\n    # Execute() narrows $dateRange down to the matched text in place, so\n    # writing back to that same Range object edits only the match (and\n    # keeps the run formatting) rather than the whole document.\n    $dateRange.Text = \"September 21, 2025\"\n}\n\n# 2. Split the mailing-address paragraph \"2940 Sanor Pl, Santa Clara CA 95051\"\n#    (the sender/recipient block near the top of the letter, not the\n#    \"PROPERTY ADDRESS:\" line further down inside the details table) into\n#    two paragraphs: \"2940 Sanor Pl\" and a new paragraph \"Santa Clara, CA 95051\".\n$addressParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $inTable = $p.Range.Information(12)  # wdWithInTable\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if (-not $inTable -and $text -eq \"2940 Sanor Pl, Santa Clara CA 95051\") {\n        $addressParagraph = $p\n        break\n    }\n}\nif ($addressParagraph -ne $null) {\n    # Setting the paragraph range's text with an embedded carriage return\n    # splits it into two paragraphs, and both inherit the original\n    # paragraph/run formatting (Arial, 22 half-points).\n    $addressParagraph.Range.Text = \"2940 Sanor Pl\" + [char]13 + \"Santa Clara, CA 95051\"\n}\n\n# 3. Remove the now-superfluous empty \"No Spacing\" paragraph that sat\n#    directly below \"...Board of Directors\".\n$paraArr = @($d.Paragraphs)\n$boardIndex = -1\nfor ($i = 0; $i -lt $paraArr.Count; $i++) {\n    if ($paraArr[$i].Range.Text.IndexOf(\"Board of Directors\") -ge 0) {\n        $boardIndex = $i\n        break\n    }\n}\nif ($boardIndex -ge 0 -and ($boardIndex + 1) -lt $paraArr.Count) {\n    $candidate = $paraArr[$boardIndex + 1]\n    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)\n    if ($candidateText -eq \"\" -and $candidate.Style.NameLocal -eq \"No Spacing\") {\n        $candidate.Range.Delete()\n    }\n}\n"}
